$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5
$ws.Range("G5").Value = 2.05
$ws.Range("I5").Value = 4
# Row 6
$ws.Range("H6").Value = 3.55
$ws.Range("I6").Value = 4.9
$ws.Range("L6").Value = 1.33
$ws.Range("M6").Value = 2.8
$ws.Range("N6").Value = 1.98
$ws.Range("O6").Value = 1.65
$ws.Range("P6").Value = 1.44
$ws.Range("Q6").Value = 2.42
$ws.Range("R6").Value = 1.93
$ws.Range("S6").Value = 1.7
$ws.Range("T6").Value = 5.9
$ws.Range("U6").Value = 7
$ws.Range("V6").Value = 8.25
$ws.Range("X6").Value = 14.5
$ws.Range("Y6").Value = 32
$ws.Range("Z6").Value = 8.75
$ws.Range("AA6").Value = 6.9
$ws.Range("AB6").Value = 18
$ws.Range("AC6").Value = 100
$ws.Range("AD6").Value = 900
$ws.Range("AE6").Value = 12
$ws.Range("AF6").Value = 28
$ws.Range("AG6").Value = 16
$ws.Range("AI6").Value = 55
$ws.Range("AJ6").Value = 60
# Row 7
$ws.Range("G7").Value = 1.34
$ws.Range("H7").Value = 4.65
$ws.Range("L7").Value = 1.2
$ws.Range("M7").Value = 3.6
$ws.Range("N7").Value = 1.62
$ws.Range("O7").Value = 2.05
$ws.Range("R7").Value = 1.93
$ws.Range("S7").Value = 1.7
$ws.Range("T7").Value = 7.1
$ws.Range("U7").Value = 6.4
$ws.Range("X7").Value = 11
$ws.Range("Y7").Value = 28
$ws.Range("Z7").Value = 13
$ws.Range("AA7").Value = 9.5
$ws.Range("AB7").Value = 21
$ws.Range("AC7").Value = 100
$ws.Range("AD7").Value = 800
$ws.Range("AE7").Value = 20
$ws.Range("AG7").Value = 24
$ws.Range("AI7").Value = 90
$ws.Range("AJ7").Value = 80
# Row 8
$ws.Range("G8").Value = 2.65
$ws.Range("H8").Value = 3.1
$ws.Range("I8").Value = 2.55
$ws.Range("L8").Value = 1.33
$ws.Range("M8").Value = 2.77
$ws.Range("N8").Value = 1.98
$ws.Range("O8").Value = 1.65
$ws.Range("P8").Value = 1.4
$ws.Range("Q8").Value = 2.52
$ws.Range("R8").Value = 1.75
$ws.Range("S8").Value = 1.87
$ws.Range("U8").Value = 13
$ws.Range("V8").Value = 9.75
$ws.Range("Z8").Value = 8.5
$ws.Range("AA8").Value = 6
$ws.Range("AB8").Value = 14
$ws.Range("AC8").Value = 70
$ws.Range("AD8").Value = 600
$ws.Range("AE8").Value = 7.7
$ws.Range("AG8").Value = 9.75
$ws.Range("AI8").Value = 23
$ws.Range("AJ8").Value = 32
# Row 9
$ws.Range("G9").Value = 1.75
$ws.Range("H9").Value = 4.25
$ws.Range("I9").Value = 3.5
$ws.Range("N9").Value = 1.34
$ws.Range("R9").Value = 1.38
$ws.Range("S9").Value = 2.57
$ws.Range("T9").Value = 13.5
$ws.Range("U9").Value = 12.5
$ws.Range("V9").Value = 9
$ws.Range("W9").Value = 17
$ws.Range("X9").Value = 12
$ws.Range("AA9").Value = 9.75
$ws.Range("AB9").Value = 12
$ws.Range("AC9").Value = 32
$ws.Range("AD9").Value = 150
$ws.Range("AE9").Value = 20
$ws.Range("AF9").Value = 26
$ws.Range("AG9").Value = 13
$ws.Range("AH9").Value = 50
$ws.Range("AI9").Value = 25
$ws.Range("AJ9").Value = 23
# Row 11
$ws.Range("G11").Value = 3.6
$ws.Range("H11").Value = 2.95
$ws.Range("I11").Value = 2.1
$ws.Range("L11").Value = 1.47
$ws.Range("M11").Value = 2.35
$ws.Range("N11").Value = 2.32
$ws.Range("O11").Value = 1.47
$ws.Range("P11").Value = 1.5
$ws.Range("Q11").Value = 2.25
$ws.Range("R11").Value = 2.02
$ws.Range("S11").Value = 1.62
$ws.Range("T11").Value = 8.25
$ws.Range("U11").Value = 18
$ws.Range("V11").Value = 13
$ws.Range("W11").Value = 55
$ws.Range("X11").Value = 40
$ws.Range("Y11").Value = 55
$ws.Range("Z11").Value = 6.6
$ws.Range("AA11").Value = 5.9
$ws.Range("AB11").Value = 18
$ws.Range("AC11").Value = 110
$ws.Range("AE11").Value = 5.7
$ws.Range("AF11").Value = 8.75
$ws.Range("AG11").Value = 9.25
$ws.Range("AH11").Value = 19.5
$ws.Range("AI11").Value = 21
$ws.Range("AJ11").Value = 40
# Row 12
$ws.Range("G12").Value = 2.3
$ws.Range("I12").Value = 3.15
$ws.Range("L12").Value = 1.47
$ws.Range("M12").Value = 2.32
$ws.Range("N12").Value = 2.37
$ws.Range("O12").Value = 1.45
$ws.Range("P12").Value = 1.52
$ws.Range("Q12").Value = 2.22
$ws.Range("R12").Value = 2.02
$ws.Range("S12").Value = 1.62
$ws.Range("T12").Value = 5.9
$ws.Range("U12").Value = 9.75
$ws.Range("W12").Value = 23
$ws.Range("X12").Value = 24
$ws.Range("Z12").Value = 6.6
$ws.Range("AA12").Value = 5.9
$ws.Range("AB12").Value = 18
$ws.Range("AC12").Value = 120
$ws.Range("AE12").Value = 7.4
$ws.Range("AF12").Value = 14.5
$ws.Range("AG12").Value = 11.75
$ws.Range("AH12").Value = 40
$ws.Range("AI12").Value = 32
$ws.Range("AJ12").Value = 50
# Row 14
$ws.Range("G14").Value = 3.5
$ws.Range("I14").Value = 2.05
$ws.Range("M14").Value = 2.32
$ws.Range("P14").Value = 1.5
$ws.Range("Q14").Value = 2.25
$ws.Range("R14").Value = 2.07
$ws.Range("S14").Value = 1.6
$ws.Range("T14").Value = 7.8
$ws.Range("U14").Value = 16.5
$ws.Range("V14").Value = 13
$ws.Range("W14").Value = 50
$ws.Range("Z14").Value = 6.8
$ws.Range("AC14").Value = 120
$ws.Range("AE14").Value = 5.6
$ws.Range("AF14").Value = 8.5
$ws.Range("AG14").Value = 9.5
$ws.Range("AH14").Value = 18
$ws.Range("AI14").Value = 21
# Row 15
$ws.Range("G15").Value = 2.05
$ws.Range("I15").Value = 3.75
$ws.Range("Q15").Value = 2.27
$ws.Range("T15").Value = 5.7
$ws.Range("U15").Value = 8.5
$ws.Range("V15").Value = 9
$ws.Range("W15").Value = 18.5
$ws.Range("X15").Value = 19.5
$ws.Range("AA15").Value = 6
$ws.Range("AB15").Value = 18.5
$ws.Range("AC15").Value = 120
$ws.Range("AE15").Value = 8
$ws.Range("AG15").Value = 13.5
$ws.Range("AH15").Value = 60
$ws.Range("AI15").Value = 45
$ws.Range("AJ15").Value = 65
